$d = $word.ActiveDocument
$om = $d.OMaths.Item(29)
$xml = $om.Range.XML()
Write-Host ("LEN=" + $xml.Length)
Write-Host $xml
